$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("N3").Value = 4.75
$ws.Range("AC3").Value = 4.75
$ws.Range("AE3").Value = 23
$ws.Range("AH3").Value = 8.5
$ws.Range("AJ3").Value = 21

# Row 4
$ws.Range("V4").Value = 1.7

# Row 5
$ws.Range("U5").Value = 2.63
$ws.Range("V5").Value = 1.44

# Row 6
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.83

# Row 7
$ws.Range("U7").Value = 1.8
$ws.Range("V7").Value = 1.91
